# Apply the header renames and per-row adaptive_filter / metric updates
# described by the commit diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header renames (row 1) ---
$ws.Range("C1").Value = "rules"
$ws.Range("E1").Value = "adaptive_filter"

# --- Per-row data (rows 2-16): column E becomes the text "RLS", and the
#     RMSE / NDEI / MAE columns (F/G/H) are refreshed with re-run values ---

$rows = @(
    @{ Row = 2;  F = 315.9341068301986;  G = 2.482342457727154;  H = 253.8352562990032 },
    @{ Row = 3;  F = 263.3272106892654;  G = 2.069002052127745;  H = 211.4722118172076 },
    @{ Row = 4;  F = 126.234718974284;   G = 0.9918454379398097; H = 102.316969203404 },
    @{ Row = 5;  F = 66.48844006468295;  G = 0.5224098131618173; H = 53.74762739848315 },
    @{ Row = 6;  F = 41.78661389888759;  G = 0.3283237979165388; H = 33.64266429758955 },
    @{ Row = 7;  F = 40.16089511457103;  G = 0.3155502775039341; H = 32.28595415376412 },
    @{ Row = 8;  F = 38.88548299743014;  G = 0.3055291699975511; H = 31.08486927236525 },
    @{ Row = 9;  F = 37.94362588552612;  G = 0.2981288550348864; H = 29.973175989316 },
    @{ Row = 10; F = 37.30799337902906;  G = 0.2931345934965549; H = 29.17586300515537 },
    @{ Row = 11; F = 36.95548023488049;  G = 0.2903648439642665; H = 28.61548169724282 },
    @{ Row = 12; F = 36.87038244280986;  G = 0.2896962176344436; H = 28.3834584565003 },
    @{ Row = 13; F = 37.00364834393676;  G = 0.2907433081428209; H = 28.41857902300262 },
    @{ Row = 14; F = 37.19231775957229;  G = 0.2922257125138006; H = 28.51983986046487 },
    @{ Row = 15; F = 37.40009775800063;  G = 0.2938582716481721; H = 28.67229858354404 },
    @{ Row = 16; F = 38.9818814560824;   G = 0.3062865873880722; H = 29.9753927343908 }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 5).Value = "RLS"
    $ws.Cells.Item($row, 6).Value = $r.F
    $ws.Cells.Item($row, 7).Value = $r.G
    $ws.Cells.Item($row, 8).Value = $r.H
}
